# Completes the test-acceptance plan: adds "produit", "panier" and
# "confirmation" test rows, switches the "Resultat observe" column to a
# centered "OK", clears stale placeholder content and re-selects E2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# --- Row heights for the three new, taller content rows ---
$ws.Rows.Item(3).RowHeight = 129.6
$ws.Rows.Item(4).RowHeight = 259.2
$ws.Rows.Item(5).RowHeight = 108

# --- Row 3: new "produit" test case (cell text set in shared-string order) ---
$ws.Cells.Item(3,2).Value = "Une page “produit” qui affiche (de manière dynamique) les détails du produit sur`nlequel l'utilisateur a cliqué depuis la page d’accueil. Depuis cette page, l’utilisateur`npeut sélectionner une quantité, une couleur, et ajouter le produit à son panier."
$ws.Cells.Item(3,4).Value = "Affichage d'un seul produit"
$ws.Cells.Item(2,5).Value = "OK"
$ws.Cells.Item(3,3).Value = "Ouvrir sur la page produit du site web dans un navigateur"

# --- Row 4: new "panier" test case ---
$ws.Cells.Item(4,2).Value = "Une page “panier”. Celle-ci contient plusieurs parties :`n○ Un résumé des produits dans le panier, le prix total et la possibilité de`nmodifier la quantité d’un produit sélectionné ou bien de supprimer celui-ci.`n○ Un formulaire permettant de passer une commande. Les données du`nformulaire doivent être correctes et bien formatées avant d'être renvoyées au`nback-end. Par exemple, pas de chiffre dans un champ prénom."
$ws.Cells.Item(4,3).Value = "Ouvrir sur la page panier du site web dans un navigateur"
$ws.Cells.Item(4,4).Value = "Affichage des produits ajouté au panier"

# --- Row 5: new "confirmation" test case ---
$ws.Cells.Item(5,2).Value = "Une page “confirmation” :`n○ Un message de confirmation de commande, remerciant l'utilisateur pour sa`ncommande, et indiquant l'identifiant de commande envoyé par l’API."
$ws.Cells.Item(5,3).Value = "Ouvrir sur la page confirmation du site web dans un navigateur"
$ws.Cells.Item(5,4).Value = "Confirmation et affichage de son numéro de commande."

# --- "Resultat observe" column: centered "OK" for every new/kept row ---
$ws.Cells.Item(3,5).Value = "OK"
$ws.Cells.Item(4,5).Value = "OK"
$ws.Cells.Item(5,5).Value = "OK"
$ws.Cells.Item(2,5).HorizontalAlignment = $xlCenter
$ws.Cells.Item(3,5).HorizontalAlignment = $xlCenter
$ws.Cells.Item(4,5).HorizontalAlignment = $xlCenter
$ws.Cells.Item(5,5).HorizontalAlignment = $xlCenter

# --- Rows 6-8: clear stale placeholder numbering/text, center "Resultat observe" ---
$ws.Cells.Item(6,1).ClearContents()
$ws.Cells.Item(6,5).HorizontalAlignment = $xlCenter
$ws.Cells.Item(7,1).ClearContents()
$ws.Cells.Item(7,5).HorizontalAlignment = $xlCenter
$ws.Cells.Item(8,5).HorizontalAlignment = $xlCenter

# --- Selection moves to E2 ---
$ws.Range("E2").Select()

# --- Page setup (A4 portrait, as set when the sheet was last printed) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

